# Applies Thai translation edits to "Email 7 [TEMPLATE] Partner email" document.
$d = $word.ActiveDocument

function ReplaceInPara($paraIndex, $find, $replace) {
    $rng = $d.Paragraphs($paraIndex).Range
    $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# Paragraph 1: language switcher line
ReplaceInPara 1 "English" "ภาษาอังกฤษ"
ReplaceInPara 1 " / Portuguese / French / Thai / Vietnamese / Spanish" " / ภาษาโปรตุเกส / ภาษาฝรั่งเศส /ภาษาไทย / ภาษาเวียดนาม / ภาษาสเปน"

# Paragraph 3: "English" heading under bookmark
ReplaceInPara 3 "English" "ภาษาอังกฤษ"

# Paragraph 5 (table cell): "Brief" label
ReplaceInPara 5 "Brief" "บทย่อ"

# Paragraph 6 (table cell): brief description
ReplaceInPara 6 "An email sent to the confirmed attendees of the event. It will be sent via customer.io" "An email sent to the confirmed attendees of the event. โดยมันจะถูกส่งผ่านทาง customer.io"

# Paragraph 8 (table cell): "Target audience" label
ReplaceInPara 8 "Target audience" "กลุ่มเป้าหมาย"

# Paragraph 13: heading "Travel checklist..."
ReplaceInPara 13 "Travel checklist: here's what you need" "รายการตรวจสอบสิ่งจำเป็นในการเดินทาง: นี่คือสิ่งที่คุณต้องมี"

# Paragraph 15: "Hi [PARTNER NAME]," -> "สวัสดี [PARTNER NAME]" (comma run removed)
ReplaceInPara 15 "Hi " "สวัสดี "
ReplaceInPara 15 "," ""

# Paragraph 16: single-day event sentence
ReplaceInPara 16 "You are all set to attend " "คุณพร้อมที่จะเข้าร่วมงาน "
ReplaceInPara 16 ", happening on " " ที่จะจัดขึ้นในวันที่ "
ReplaceInPara 16 " at " " ที่ "
ReplaceInPara 16 "!" " แล้ว!"

# Paragraph 17: multi-day event sentence
ReplaceInPara 17 "You are all set to attend " "คุณพร้อมที่จะเข้าร่วมงาน "
ReplaceInPara 17 ", happening from " " ที่จะจัดขึ้นตั้งแต่วันที่ "
ReplaceInPara 17 " to " " ถึงวันที่ "
ReplaceInPara 17 " at " " ที่ "
ReplaceInPara 17 "! " " แล้ว! "

# Paragraph 18: checklist intro
ReplaceInPara 18 "Here’s a checklist of the necessary items for your trip: " "นี่คือรายการตรวจสอบสิ่งจำเป็นสำหรับทริปการเดินทางของคุณ: "

# Paragraph 19: Passport
ReplaceInPara 19 "Passport " "หนังสือเดินทาง "

# Paragraph 20: Visa (if applicable)
ReplaceInPara 20 "Visa " "วีซ่า "
ReplaceInPara 20 "(if applicable) " "(ถ้ามี) "

# Paragraph 22: yellow fever vaccination certificate bullet
ReplaceInPara 22 "Valid yellow fever vaccination certificate " "ใบรับรองการฉีดวัคซีนไข้เหลืองที่ยังไม่หมดอายุ "
ReplaceInPara 22 "For travellers from yellow fever endemic countries, follow the requirements set by your country. Vaccination should be done no less than 14 days prior to the journey. " " สำหรับผู้เดินทางจากประเทศที่มีไข้เหลืองเป็นโรคเฉพาะถิ่น โปรดปฏิบัติตามเงื่อนไขที่กำหนดโดยประเทศของคุณ การฉีดวัคซีนควรจะทำไม่น้อยกว่า 14 วันก่อนการเดินทาง "
ReplaceInPara 22 "(As you’re travelling on " "(ในกรณีที่คุณจะเดินทางในวันที่ "
ReplaceInPara 22 ", you should have received your vaccination before or on" " คุณควรได้รับการฉีดวัคซีนก่อนหรือในวันที่"
ReplaceInPara 22 ".)" ")"

# Paragraph 23: digital/printed copy of travel itinerary
ReplaceInPara 23 "A digital or printed copy of the travel itinerary" "สำเนากำหนดการเดินทาง (travel itinerary) ในรูปแบบดิจิทัลหรือพิมพ์ออกมา"

# Paragraph 24: smart casual attire
ReplaceInPara 24 "Smart casual attire for the conference" "เสื้อผ้าแบบสมาร์ทแคชชวลสำหรับใส่เข้าร่วมงานประชุม"

# Paragraph 25: black tie attire
ReplaceInPara 25 "Black tie attire for the Gala dinner" "เสื้อผ้าแบล็กไท (Black tie attire) สำหรับงานเลี้ยงกาลาดินเนอร์"

# Paragraph 29: contact via live chat or WhatsApp
ReplaceInPara 29 "If you have any questions, please contact us via " "หากคุณมีคำถามใดๆ กรุณาติดต่อเราผ่านทาง "
ReplaceInPara 29 "live chat" "แชทสด"
ReplaceInPara 29 " or " " หรือทาง "
ReplaceInPara 29 ". " " "

# Paragraph 30: contact country manager
ReplaceInPara 30 "If you have any questions, please contact your country manager, " "หากคุณมีคำถามใดๆ โปรดติดต่อผู้จัดการประจำประเทศของคุณซึ่งได้แก่ "
ReplaceInPara 30 ", at " " ที่ "
ReplaceInPara 30 " or " " หรือ "
ReplaceInPara 30 " (WhatsApp). " " (WhatsApp) "

# Comment (w:id=5): "choose either one"
foreach ($c in $d.Comments) {
    if ($c.Range.Text -eq "choose either one") {
        $c.Range.Text = "เลือกอย่างใดอย่างหนึ่ง"
    }
}
